$d = $word.ActiveDocument

$d.Content.Find.Execute("activitate_id", $true, $false, $false, $false, $false,
                         $true, 1, $false, "grup_id", 2)
